## Estado de Cuenta - add a second worker row, update totals (NIT-9000950342)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new worker data row (row 17) below the existing one (row 16). ---
# Copy row 16 so the insert carries its formatting down and pushes the
# signature block (old rows 21-22) down to rows 22-23.
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()

# Re-apply row 16's cell formatting (borders/fonts/number formats) onto the
# freshly inserted row 17 so it keeps the same bordered-table look.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill in the new worker's data in row 17. ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002413552"
$ws.Range("D17").Value = "PABLO JOSE OSPINO PEREZ"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 21000
$ws.Range("G17").Value = 1750000

# --- 3. Update the account summary numbers. ---
# Total "Valor Mora" at the top now reflects both workers combined.
$ws.Range("E11").Value = 36600

# Worker / period counts go from 1 to 2 now that a second worker was added.
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
